# For each plate-reader data sheet (WT, N121E, F193D, F159G), the "Part of
# Plate"/"B1-G12" summary row is removed, and the B-G data block (rows 33-38)
# is bracketed by two new, empty placeholder rows labelled "A" (row 32) and
# "H" (row 39) -- turning the table into a full A-H plate-row listing. The
# selection on every sheet is updated to A32:A39 (the new label column) and
# the previous topLeftCell scroll position is cleared.

$wb = $excel.ActiveWorkbook

$sheetNames = @("WT", "N121E", "F193D", "F159G")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the "Part of Plate" / "B1-G12" row entirely; everything below
    # shifts up by one (old row 29 "Start Time" becomes row 28, etc.).
    $ws.Rows.Item(28).Delete()

    # Make room for a new label row ("A") right after the "<>" header row
    # (old row 32, now row 31) and above the "B" data row (old row 33, now
    # row 32 post-delete).
    $ws.Rows.Item(32).Insert()

    # Give the new row the same look as the other plate-row labels (white
    # text on grey fill, style used by A31:A38) by copying the format from
    # the row right below it.
    $ws.Range("A38").Copy() | Out-Null
    $ws.Range("A32").PasteSpecial(-4122) | Out-Null
    $ws.Range("A32").Value = "A"

    # Row 39 is already blank (it was empty, unused space between the "G"
    # data row and "End Time:"), so just populate it directly -- no insert
    # needed here, which keeps "End Time:" pinned at row 42.
    $ws.Range("A38").Copy() | Out-Null
    $ws.Range("A39").PasteSpecial(-4122) | Out-Null
    $ws.Range("A39").Value = "H"

    $excel.CutCopyMode = 0

    # Refresh the view: select the full label column A32:A39 and drop any
    # stale scroll (topLeftCell) position.
    $ws.Activate()
    $ws.Range("A32:A39").Select()
}

# F159G (the 4th sheet) is the tab that was active/selected last.
$wb.Worksheets.Item("F159G").Activate()
$wb.Worksheets.Item("F159G").Range("A32:A39").Select()
